# Reverses the order of the "Periodo Mora" (column E) and "Valor Mora"
# (column F) data for rows 16-22 on Hoja1 - old EC periods removed, new
# ones added (commit: "Elimna EC anteriores y se agregan nuevos, se
# modifica base de datos").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Column E (Periodo Mora) - reverse the 7-row block (rows 16..22)
$ws.Range("E16").Value = "2506"
$ws.Range("E17").Value = "2505"
$ws.Range("E18").Value = "2504"
$ws.Range("E19").Value = "2503"
$ws.Range("E20").Value = "2502"
$ws.Range("E21").Value = "2501"
$ws.Range("E22").Value = "2412"

# Column F (Valor Mora) follows the same reversal - only the first and
# last rows actually differ in value (52000 <-> 45066), the middle ones
# stay 52000.
$ws.Range("F16").Value = 45066
$ws.Range("F22").Value = 52000
